$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 48 - more precise timestamp value
$ws.Cells.Item(48, 1).Value = 45913.92153459491

# New row 49
$ws.Cells.Item(49, 1).Value = 45914.3330212963
$ws.Cells.Item(49, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(49, 2).Value = 660
$ws.Cells.Item(49, 3).Value = 264
$ws.Cells.Item(49, 4).Value = "0:0"
$ws.Cells.Item(49, 5).Value = "Decrease Unit load to 264.0 MW/Giảm tải xuống 264.0 MW"

# New row 50
$ws.Cells.Item(50, 1).Value = 45914.34451471065
$ws.Cells.Item(50, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(50, 2).Value = 660
$ws.Cells.Item(50, 3).Value = 264
$ws.Cells.Item(50, 4).Value = "00:00"
$ws.Cells.Item(50, 5).Value = "Decrease Unit load to 264.0 MW/Giảm tải xuống 264.0 MW"

# New row 51
$ws.Cells.Item(51, 1).Value = 45914.34880066659
$ws.Cells.Item(51, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(51, 2).Value = 660
$ws.Cells.Item(51, 3).Value = 462
$ws.Cells.Item(51, 4).Value = "00:00"
$ws.Cells.Item(51, 5).Value = "Decrease Unit load to 462.0 MW/Giảm tải xuống 462.0 MW"

# New row 52
$ws.Cells.Item(52, 1).Value = 45914.34907038009
$ws.Cells.Item(52, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(52, 2).Value = 660
$ws.Cells.Item(52, 3).Value = 264
$ws.Cells.Item(52, 4).Value = "00:00"
$ws.Cells.Item(52, 5).Value = "Decrease Unit load to 264.0 MW/Giảm tải xuống 264.0 MW"
